# Prefix each protocol/step sheet's own name onto the "Name" column (A)
# values of its data rows (everything below the header row), e.g.
# "Step4 Seed" -> "free1 Step4 Seed" on sheet "free1".
#
# This mirrors the commit:
#   fix: unique command names in XLSX - prefix protocol name to each step

$wb = $excel.ActiveWorkbook

# Sheets whose column-A "Name" values get the sheet name prefixed.
# (The first five overview/meta sheets are left untouched.)
$targetSheets = @(
    "price1", "price2",
    "discount1", "discount2",
    "free1", "free2",
    "nomoney1", "nomoney2",
    "noppv1", "noppv2",
    "card1", "card2",
    "nosex1", "nosex2",
    "offtopic1", "offtopic2",
    "real1", "real2",
    "voice1", "voice2",
    "customyes1", "customyes2",
    "customno1", "customno2",
    "done1", "done2",
    "cumcontrol", "dickpic", "boosters"
)

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $usedRange = $ws.UsedRange
    $lastRow = $usedRange.Rows.Count

    # Row 1 is the header ("Name", "Text", "Note", "*Guidelines"); data
    # starts at row 2.
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $current = $cell.Text
        if ([string]::IsNullOrEmpty($current)) {
            continue
        }
        $prefix = "$sheetName "
        if ($current.StartsWith($prefix)) {
            continue
        }
        $cell.Value = "$sheetName $current"
    }
}
